$wb = $excel.ActiveWorkbook

# Overview sheet: G2 holds "Latest HO Xliff Generate Date" for the .md file row.
# This shared string is also used by de-de!H2 (same timestamp text), so updating
# either cell updates the shared string used by both.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 06:57:44"

# zh-cn sheet: H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 06:57:39"
$wsZhCn.Range("K2").Value = "2016-08-17 06:57:56"

# de-de sheet: H2 = Correspond Handoff Datetime (same value as Overview!G2),
# K2 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 06:57:44"
$wsDeDe.Range("K2").Value = "2016-08-17 06:58:09"
